$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last row (row 9); table shrinks from 8 data rows to 7 data rows
$ws.Rows.Item(9).Delete()

# Update remaining data rows (2-8) with the new test results
$ws.Range("A2").Value = "ADF Test (Drift)"
$ws.Range("B2").Value = "Unit Root"
$ws.Range("C2").Value = -1.67
$ws.Range("D2").Value = -2.89
$ws.Range("E2").Value = "Not Stationary"
$ws.Range("F2").Value = $true
$ws.Range("A3").Value = "ADF Test (Differenced)"
$ws.Range("B3").Value = "Unit Root"
$ws.Range("C3").Value = -14.12
$ws.Range("D3").Value = -2.89
$ws.Range("E3").Value = "Not stationary"
$ws.Range("F3").Value = $false
$ws.Range("A4").Value = "Phillips-Perron Test"
$ws.Range("B4").Value = "Unit Root"
$ws.Range("C4").Value = -18.26
$ws.Range("D4").Value = -2.89
$ws.Range("E4").Value = "Not Stationary"
$ws.Range("F4").Value = $false
$ws.Range("A5").Value = "ERS Test (DF-GLS)"
$ws.Range("B5").Value = "Unit Root"
$ws.Range("C5").Value = -4.97
$ws.Range("D5").Value = -1.94
$ws.Range("E5").Value = "Not stationary"
$ws.Range("F5").Value = $false
$ws.Range("A6").Value = "ERS Test (P-test)"
$ws.Range("B6").Value = "Unit Root"
$ws.Range("C6").Value = 0.52
$ws.Range("D6").Value = 3.11
$ws.Range("E6").Value = "Not Stationary"
$ws.Range("F6").Value = $false
$ws.Range("A7").Value = "KPSS Test (Tau)"
$ws.Range("B7").Value = "Unit Root"
$ws.Range("C7").Value = 0.07
$ws.Range("D7").Value = 0.15
$ws.Range("E7").Value = "Stationary"
$ws.Range("F7").Value = $true
$ws.Range("A8").Value = "KPSS Test (Mu)"
$ws.Range("B8").Value = "Unit Root"
$ws.Range("C8").Value = 0.07
$ws.Range("D8").Value = 0.46
$ws.Range("E8").Value = "Stationary"
$ws.Range("F8").Value = $true
